$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove previously-present descriptive stats in rows 3-5 for these columns
$ws.Range("G3").ClearContents()
$ws.Range("H3").ClearContents()
$ws.Range("I3").ClearContents()
$ws.Range("J3").ClearContents()
$ws.Range("K3").ClearContents()
$ws.Range("L3").ClearContents()
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()
$ws.Range("O3").ClearContents()
$ws.Range("P3").ClearContents()
$ws.Range("Q3").ClearContents()
$ws.Range("R3").ClearContents()
$ws.Range("S3").ClearContents()
$ws.Range("T3").ClearContents()
$ws.Range("U3").ClearContents()
$ws.Range("V3").ClearContents()
$ws.Range("W3").ClearContents()
$ws.Range("Y3").ClearContents()
$ws.Range("AC3").ClearContents()
$ws.Range("AG3").ClearContents()
$ws.Range("G4").ClearContents()
$ws.Range("H4").ClearContents()
$ws.Range("I4").ClearContents()
$ws.Range("J4").ClearContents()
$ws.Range("K4").ClearContents()
$ws.Range("L4").ClearContents()
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("O4").ClearContents()
$ws.Range("P4").ClearContents()
$ws.Range("Q4").ClearContents()
$ws.Range("R4").ClearContents()
$ws.Range("S4").ClearContents()
$ws.Range("T4").ClearContents()
$ws.Range("U4").ClearContents()
$ws.Range("V4").ClearContents()
$ws.Range("W4").ClearContents()
$ws.Range("Y4").ClearContents()
$ws.Range("AC4").ClearContents()
$ws.Range("AG4").ClearContents()
$ws.Range("G5").ClearContents()
$ws.Range("H5").ClearContents()
$ws.Range("I5").ClearContents()
$ws.Range("J5").ClearContents()
$ws.Range("K5").ClearContents()
$ws.Range("L5").ClearContents()
$ws.Range("M5").ClearContents()
$ws.Range("N5").ClearContents()
$ws.Range("O5").ClearContents()
$ws.Range("P5").ClearContents()
$ws.Range("Q5").ClearContents()
$ws.Range("R5").ClearContents()
$ws.Range("S5").ClearContents()
$ws.Range("T5").ClearContents()
$ws.Range("U5").ClearContents()
$ws.Range("V5").ClearContents()
$ws.Range("W5").ClearContents()
$ws.Range("Y5").ClearContents()
$ws.Range("AC5").ClearContents()
$ws.Range("AG5").ClearContents()

# Add newly-computed descriptive stats for rows 6-12
$ws.Range("G6").Value = 51.48472391399029
$ws.Range("H6").Value = 50.60755926316322
$ws.Range("I6").Value = 48.15822712641475
$ws.Range("J6").Value = 56.38270831546305
$ws.Range("K6").Value = 2653.012561327719
$ws.Range("L6").Value = 2873.343491797922
$ws.Range("M6").Value = 2732.232641521527
$ws.Range("N6").Value = 165106.2475409963
$ws.Range("O6").Value = 19.8904978311913
$ws.Range("P6").Value = 22.03197717961866
$ws.Range("Q6").Value = 95.79049408224471
$ws.Range("R6").Value = 0.925780582370485
$ws.Range("S6").Value = 0.8439787173411903
$ws.Range("T6").Value = 3.339555440648116
$ws.Range("U6").Value = 50.6098520466497
$ws.Range("V6").Value = 9408702754.195225
$ws.Range("W6").Value = 53.98617537798903
$ws.Range("Y6").Value = 0.1508458506960731
$ws.Range("AC6").Value = 18.5595401452801
$ws.Range("AG6").Value = 8.181976829350004
$ws.Range("G7").Value = 13.61498356546453
$ws.Range("H7").Value = 20.27993945625253
$ws.Range("I7").Value = 18.08384025949885
$ws.Range("J7").Value = 19.87157400581872
$ws.Range("K7").Value = 1723.505173154446
$ws.Range("L7").Value = 1603.510961402266
$ws.Range("M7").Value = 1499.288314395472
$ws.Range("N7").Value = 121628.1656343791
$ws.Range("O7").Value = 41.28561996785294
$ws.Range("P7").Value = 66.45499500471675
$ws.Range("Q7").Value = 2385.801069781274
$ws.Range("R7").Value = 0.6094216696086281
$ws.Range("S7").Value = 0.5427660906465637
$ws.Range("T7").Value = 18.37286716499127
$ws.Range("U7").Value = 25.87312699073571
$ws.Range("V7").Value = 31025589787.37921
$ws.Range("W7").Value = 29.7272498150084
$ws.Range("Y7").Value = 0.487463762369141
$ws.Range("AC7").Value = 354.6696458411183
$ws.Range("AG7").Value = 124.5591951743139
$ws.Range("G8").Value = 11.9320378601968
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 2.68130501339456
$ws.Range("J8").Value = 6.78899075483514
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 124.8114555427984
$ws.Range("N8").Value = 1698.806320386001
$ws.Range("O8").Value = 0.028
$ws.Range("P8").Value = 0.096212113990482
$ws.Range("Q8").Value = 0.06748829419609099
$ws.Range("R8").Value = -0.465974128803762
$ws.Range("S8").Value = -0.13823804827088
$ws.Range("T8").Value = -275.8579110609
$ws.Range("U8").Value = 0
$ws.Range("V8").Value = 8139905.5198973
$ws.Range("W8").Value = 0
$ws.Range("Y8").Value = [double]"1.118963486454647e-05"
$ws.Range("AC8").Value = 0.030235587
$ws.Range("AG8").Value = 0.006938358022568212
$ws.Range("G9").Value = 42.5218530970591
$ws.Range("H9").Value = 36.8064951146347
$ws.Range("I9").Value = 34.5615397966496
$ws.Range("J9").Value = 41.6182456140351
$ws.Range("K9").Value = 1354.170566090234
$ws.Range("L9").Value = 1584.988432018132
$ws.Range("M9").Value = 1617.909969968546
$ws.Range("N9").Value = 76884.10236875113
$ws.Range("O9").Value = 2.35
$ws.Range("P9").Value = 3.625497746698368
$ws.Range("Q9").Value = 5.842975093274302
$ws.Range("R9").Value = 0.589173640877774
$ws.Range("S9").Value = 0.471751656461634
$ws.Range("T9").Value = 1.1752825705
$ws.Range("U9").Value = 34.9625798581393
$ws.Range("V9").Value = 1171052432.04107
$ws.Range("W9").Value = 36.0482416946772
$ws.Range("Y9").Value = 0.04605633531471059
$ws.Range("AC9").Value = 0.5435058125007795
$ws.Range("AG9").Value = 0.4444574111924894
$ws.Range("G10").Value = 51.2230762947906
$ws.Range("H10").Value = 51.2518538306358
$ws.Range("I10").Value = 46.6322699828008
$ws.Range("J10").Value = 57.1733610421597
$ws.Range("K10").Value = 2361.082820166388
$ws.Range("L10").Value = 2722.742661527486
$ws.Range("M10").Value = 2491.627453346198
$ws.Range("N10").Value = 134399.2891732987
$ws.Range("O10").Value = 7.156548
$ws.Range("P10").Value = 9.1736587605
$ws.Range("Q10").Value = 10.63272990335615
$ws.Range("R10").Value = 0.902749139469404
$ws.Range("S10").Value = 0.759634935815751
$ws.Range("T10").Value = 3.74589001185
$ws.Range("U10").Value = 45.6255990665043
$ws.Range("V10").Value = 3361633598.04034
$ws.Range("W10").Value = 49.5824211135239
$ws.Range("Y10").Value = 0.09404923877341792
$ws.Range("AC10").Value = 1.462185443
$ws.Range("AG10").Value = 1.018089826734382
$ws.Range("G11").Value = 60.6964300536413
$ws.Range("H11").Value = 65.1536086515201
$ws.Range("I11").Value = 61.1923642947161
$ws.Range("J11").Value = 72.7691228070175
$ws.Range("K11").Value = 3676.945963857558
$ws.Range("L11").Value = 3901.910634991028
$ws.Range("M11").Value = 3569.112877049816
$ws.Range("N11").Value = 223609.0850257536
$ws.Range("O11").Value = 19.06
$ws.Range("P11").Value = 19.48876044822527
$ws.Range("Q11").Value = 21.71256451800405
$ws.Range("R11").Value = 1.19078344892068
$ws.Range("S11").Value = 1.17043133124314
$ws.Range("T11").Value = 7.54433891875
$ws.Range("U11").Value = 61.8064999988478
$ws.Range("V11").Value = 8475878901.60462
$ws.Range("W11").Value = 65.57001639811919
$ws.Range("Y11").Value = 0.1711457026153066
$ws.Range("AC11").Value = 3.7642770419078
$ws.Range("AG11").Value = 2.668532761542454
$ws.Range("G12").Value = 91.0727090181734
$ws.Range("H12").Value = 98.2761721811863
$ws.Range("I12").Value = 95.22303729846359
$ws.Range("J12").Value = 97.8127413127413
$ws.Range("K12").Value = 9195.123247177175
$ws.Range("L12").Value = 7981.226700483277
$ws.Range("M12").Value = 8110.709913728177
$ws.Range("N12").Value = 755378.7537650267
$ws.Range("O12").Value = 519.92648879
$ws.Range("P12").Value = 1414.073800429
$ws.Range("Q12").Value = 89368.42105263199
$ws.Range("R12").Value = 3.29542616368157
$ws.Range("S12").Value = 2.69388553673174
$ws.Range("T12").Value = 377.817151945
$ws.Range("U12").Value = 230.045673051157
$ws.Range("V12").Value = 697982117436.631
$ws.Range("W12").Value = 319.428684596163
$ws.Range("Y12").Value = 14.81738443550587
$ws.Range("AC12").Value = 13016.1328463612
$ws.Range("AG12").Value = 4392.838036620215
